$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 describe two different observations recorded at the same
# place/date by the same reporter. The species-specific fields (id,
# taxon id/name, coordinates, amount, substrate note) were swapped
# between the two rows; everything else (status, location text,
# accuracy, dates, observer, ...) stays put.

function Set-TextCell($cell, $value) {
    if ($value -eq $null) {
        # Source side had no cell at all -> remove this one entirely.
        $cell.ClearContents()
    } else {
        # Force text type even for a digit-only string (e.g. "1") or an
        # empty one, so the cell doesn't silently become a number (or a
        # fully-cleared cell) on write-back. A leading apostrophe is how
        # Excel itself marks "treat as text" on entry.
        $cell.Value = "'" + $value
        $cell.ClearFormats()
    }
}

# --- Plain numeric columns: straightforward value swap ---
$numericCols = @("A", "B", "E", "Q", "R")
foreach ($col in $numericCols) {
    $c5 = $ws.Range($col + "5")
    $c6 = $ws.Range($col + "6")
    $v5 = $c5.Value2
    $v6 = $c6.Value2
    $c5.Value2 = $v6
    $c6.Value2 = $v5
}

# --- Ordinary text columns: swap, no numeric-lookalike concerns ---
$textCols = @("F", "G", "H")
foreach ($col in $textCols) {
    $c5 = $ws.Range($col + "5")
    $c6 = $ws.Range($col + "6")
    $v5 = $c5.Value2
    $v6 = $c6.Value2
    $c5.Value2 = $v6
    $c6.Value2 = $v5
}

# --- "Antal" (I): text column whose content looks numeric ("1") or is
#     blank, so it needs the forced-text helper in both directions.
$i5 = $ws.Range("I5")
$i6 = $ws.Range("I6")
$vi5 = $i5.Value2
$vi6 = $i6.Value2
Set-TextCell $i5 $vi6
Set-TextCell $i6 $vi5

# --- Substrate description (AO): only one of the two rows has it. The
#     row losing it should end up with no cell there at all.
$ao5 = $ws.Range("AO5")
$ao6 = $ws.Range("AO6")
$vao5 = $ao5.Value2
$vao6 = $ao6.Value2
Set-TextCell $ao5 $vao6
Set-TextCell $ao6 $vao5
